$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (IDPT) - updated values
$ws.Range("B2").Value = 0.9571788124371639
$ws.Range("C2").Value = 0.3802022503362134
$ws.Range("G2").Value = 0.3186155176884216
$ws.Range("H2").Value = 0.8089986761247983
$ws.Range("I2").Value = 0.5205558643947931
$ws.Range("J2").Value = 0.4710534615481041
$ws.Range("K2").Value = 0.7023509530237569
$ws.Range("L2").Value = 0.3527039604500516
$ws.Range("P2").Value = -0.1340323814162643
$ws.Range("Q2").Value = -0.1639966123942181
$ws.Range("R2").Value = 2.168645455033084

# Row 3 - rename method from SPCT to GDPT
$ws.Range("A3").Value = "GDPT"
$ws.Range("C3").Value = 1.511179703834972
$ws.Range("G3").Value = 0.5389727935843093
$ws.Range("H3").Value = 1.103888480179478
$ws.Range("I3").Value = 0.3382432830856023
$ws.Range("J3").Value = 0.4237430169769091
$ws.Range("K3").Value = 0.5475578821850129
$ws.Range("L3").Value = 0.9217035104729026
$ws.Range("P3").Value = -0.472411337288048
$ws.Range("Q3").Value = -0.06789915326352265
$ws.Range("R3").Value = 2.185571852211659
